$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.093.23"
$ws.Range("E2").Value = "  -3.66%  "
$ws.Range("D3").Value = "1.645.45"
$ws.Range("E3").Value = "  -3.47%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'307.55"
$ws.Range("E5").Value = "  -2.77%  "
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "'0.3894"
$ws.Range("E7").Value = "  -2.76%  "
$ws.Range("D8").Value = "'0.3851"
$ws.Range("E8").Value = "  -4.72%  "
$ws.Range("D9").Value = "'1.002"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'1.345"
$ws.Range("E10").Value = "  -8.71%  "
$ws.Range("D11").Value = "'48.79"
$ws.Range("E11").Value = "  -8.50%  "
$ws.Range("D12").Value = "'0.08449"
$ws.Range("E12").Value = "  -4.24%  "
$ws.Range("D13").Value = "'23.76"
$ws.Range("E13").Value = "  -8.94%  "
$ws.Range("D14").Value = "'7.112"
$ws.Range("E14").Value = "  -4.99%  "
$ws.Range("D15").Value = "'0.00001280"
$ws.Range("E15").Value = "  -5.39%  "
$ws.Range("D16").Value = "'7.472"
$ws.Range("E16").Value = "  -6.25%  "
$ws.Range("D17").Value = "1.640.08"
$ws.Range("E17").Value = "  -5.02%  "
$ws.Range("D18").Value = "'94.79"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").Value = "'0.06974"
$ws.Range("E19").Value = "  -3.11%  "
$ws.Range("D20").Value = "'20.64"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "'6.916"
$ws.Range("E21").Value = "  -5.51%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'13.60"
$ws.Range("E23").Value = "  -5.02%  "
$ws.Range("D24").Value = "24.088.72"
$ws.Range("E24").Value = "  -3.67%  "
$ws.Range("D25").Value = "'2.328"
$ws.Range("E25").Value = "  -3.25%  "
$ws.Range("D26").Value = "'2.695"
$ws.Range("E26").Value = "  -8.46%  "
$ws.Range("D27").Value = "'22.42"
$ws.Range("E27").Value = "  -4.97%  "
$ws.Range("D28").Value = "'157.88"
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("D29").Value = "'8.667"
$ws.Range("E29").Value = "  +2.86%  "
$ws.Range("D30").Value = "'141.20"
$ws.Range("E30").Value = "  -7.29%  "
$ws.Range("D31").Value = "'5.247"
$ws.Range("E31").Value = "  -13.54%  "
$ws.Range("D32").Value = "'2.459"
$ws.Range("E32").Value = "  -7.77%  "
$ws.Range("D33").Value = "1.820.58"
$ws.Range("E33").Value = "  -4.91%  "
$ws.Range("D34").Value = "'6.952"
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("D35").Value = "'0.08011"
$ws.Range("E35").Value = "  -7.09%  "
$ws.Range("D36").Value = "'0.02904"
$ws.Range("E36").Value = "  -8.38%  "
$ws.Range("D37").Value = "'0.9578"
$ws.Range("E37").Value = "  -8.91%  "
$ws.Range("D38").Value = "'0.2688"
$ws.Range("E38").Value = "  -8.01%  "
$ws.Range("E39").Value = "  -5.44%  "
$ws.Range("D40").Value = "'1.464"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("D41").Value = "'9.916"
$ws.Range("E41").Value = "  -10.46%  "
$ws.Range("D42").Value = "'0.7587"
$ws.Range("E42").Value = "  -8.16%  "
$ws.Range("D43").Value = "'13.05"
$ws.Range("E43").Value = "  -7.09%  "
$ws.Range("D44").Value = "'15.97"
$ws.Range("E44").Value = "  -6.62%  "
$ws.Range("D45").Value = "'0.6888"
$ws.Range("E45").Value = "  -6.75%  "
$ws.Range("D46").Value = "'2.474"
$ws.Range("E46").Value = "  -8.11%  "
$ws.Range("D47").Value = "'4.093"
$ws.Range("E47").Value = "  -3.76%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").Value = "'0.08333"
$ws.Range("E49").Value = "  -10.09%  "
$ws.Range("D50").Value = "'133.34"
$ws.Range("E50").Value = "  -4.76%  "
$ws.Range("D51").Value = "'1.257"
$ws.Range("E51").Value = "  -10.55%  "
